$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the commit diff.
# All target cells are plain text cells (inlineStr) in the original workbook,
# so we force a text NumberFormat before assignment (then restore the default
# "Normal" style) to avoid Excel auto-converting numeric-looking strings such
# as "1.000" or "91.80" into numbers and losing formatting / trailing zeros.
$changes = [ordered]@{
    "D2" = "29.478.94"
    "E2" = "  +0.54%  "
    "D3" = "1.877.99"
    "E3" = "  +0.46%  "
    "D4" = "0.9995"
    "D5" = "0.7141"
    "E5" = "  +0.48%  "
    "D6" = "242.11"
    "E6" = "  +0.49%  "
    "D7" = "1.000"
    "E7" = "  -0.12%  "
    "D8" = "0.3118"
    "E8" = "  +1.15%  "
    "D9" = "0.07735"
    "E9" = "  -1.82%  "
    "D10" = "25.12"
    "E10" = "  -0.88%  "
    "D11" = "0.08404"
    "E11" = "  +1.91%  "
    "D12" = "1.896.21"
    "E12" = "  +0.68%  "
    "D13" = "5.263"
    "E13" = "  +0.54%  "
    "D14" = "0.7195"
    "E14" = "  -0.34%  "
    "D15" = "91.80"
    "E15" = "  +1.20%  "
    "D16" = "29.490.55"
    "E16" = "  +0.59%  "
    "D17" = "0.000008253"
    "E17" = "  +5.56%  "
    "D18" = "6.003"
    "E18" = "  +2.91%  "
    "E19" = "  +0.29%  "
    "B20" = "Avalanche"
    "C20" = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
    "D20" = "13.26"
    "E20" = "  +0.48%  "
    "B21" = "WrappedliquidstakedEther2.0"
    "C21" = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
    "D21" = "2.131.29"
    "E21" = "  +1.22%  "
    "D22" = "0.9996"
    "E22" = "  -0.17%  "
    "D23" = "7.959"
    "E23" = "  -0.62%  "
    "E24" = "  -0.21%  "
    "D25" = "0.1630"
    "E25" = "  +2.33%  "
    "D26" = "163.80"
    "E26" = "  +0.88%  "
    "D27" = "9.055"
    "E27" = "  +0.71%  "
    "D28" = "18.68"
    "E28" = "  +2.50%  "
    "D29" = "1.511"
    "E29" = "  +1.15%  "
    "D30" = "4.433"
    "E30" = "  +1.27%  "
    "E31" = "  -4.08%  "
    "D32" = "4.325"
    "E32" = "  +5.59%  "
    "D33" = "0.05243"
    "E33" = "  +1.15%  "
    "D34" = "1.939"
    "E34" = "  +0.04%  "
    "D35" = "0.7750"
    "E35" = "  +7.61%  "
    "D36" = "1.179"
    "E36" = "  -0.65%  "
    "D37" = "2.680"
    "E37" = "  +0.35%  "
    "D38" = "0.01870"
    "E38" = "  +0.77%  "
    "D39" = "2.726"
    "D40" = "1.177.21"
    "E40" = "  +0.40%  "
    "D41" = "6.437"
    "E41" = "  +5.41%  "
    "D42" = "73.77"
    "E42" = "  +1.65%  "
    "D43" = "0.8923"
    "E43" = "  -1.14%  "
    "D44" = "104.42"
    "E44" = "  +2.44%  "
    "D45" = "0.9993"
    "E45" = "  -0.17%  "
    "D46" = "2.030.01"
    "E46" = "  +0.94%  "
    "E47" = "  +1.13%  "
    "D48" = "0.5205"
    "E48" = "  -1.61%  "
    "D49" = "9.446"
    "E49" = "  +2.11%  "
    "D50" = "0.4327"
    "E50" = "  +1.16%  "
    "D51" = "7.088"
    "E51" = "  +1.08%  "
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.Style = "Normal"
}
